$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: Cristal 8MHz
$ws.Range("A14").Value = "YG8M000000S418"
$ws.Range("E14").Value = "Cristal 8MHz"
$ws.Range("F14").Value = "https://www.digikey.ca/en/products/detail/nextgen-components/YG8M000000S418/17289334"

# Row 15: Connecteur SMA 90deg
$ws.Range("F15").Value = "https://www.digikey.ca/en/products/detail/adam-tech/RF2-03E-T-00-50-G/9831386"
$ws.Range("A15").Value = "RF2-03E-T-00-50-G"
$ws.Range("E15").Value = "Connecteur SMA 90deg"

# Remaining numeric / common cells
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1.85
$ws.Range("G14").Value = "."

$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 4.19
$ws.Range("G15").Value = "."

# Move selection to match the new edited cell
$null = $ws.Range("E15").Select()
